# Weekly update for "Hortaliza, Vega Modelo de Temuco - Cebolla":
# 5 new daily records are inserted at the top of the most-recent block
# (rows 811-815), pushing the previously-last 12 rows (811-822) down to
# become rows 816-827. The dimension grows from A1:R822 to A1:R827.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows at 811..815 (shifts existing 811:822 down to 816:827)
$ws.Range("A811:A815").EntireRow.Insert()

# Helper data for the 5 new rows. Columns:
# A Mercado ID, B Mercado, C Región, D Fecha, E Codreg, F Categoría ID,
# G Categoría, H Variedad, I Calidad, J Volumen, K Precio mínimo,
# L Precio máximo, M Precio promedio ponderado, N Unidad de comercialización,
# O Origen, P Precio $/Kg, Q Kg o Unidades, R Clasificación
$newRows = @(
    @{ Row=811; D=44595; H="Morada(o)";        I="1a (guarda)";  J=250;  K=10000; L=12000; M=11200; N="$/malla 18 kilos";                               O="Región de O'Higgins"; P=622; Q=18 },
    @{ Row=812; D=44595; H="Sin especificar";   I="1a nueva(o)";  J=600;  K=5000;  L=5000;  M=5000;  N="$/malla 18 kilos";                               O="Región de O'Higgins"; P=278; Q=18 },
    @{ Row=813; D=44595; H="Sin especificar";   I="1a nueva(o)";  J=1000; K=4500;  L=5000;  M=4750;  N="$/malla 18 kilos";                               O="Región del Maule";    P=264; Q=18 },
    @{ Row=814; D=44595; H="Sin especificar";   I="1a nueva(o)";  J=1500; K=1000;  L=1000;  M=1000;  N="$/paquete 10 unidades (volumen en unidades)";    O="Región del Maule";    P=100; Q=10 },
    @{ Row=815; D=44595; H="Sin especificar";   I="Primera";      J=800;  K=4500;  L=4500;  M=4500;  N="$/malla 18 kilos";                               O="Perú";                P=250; Q=18 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = 10
    $ws.Cells.Item($row, 2).Value  = "Vega Modelo de Temuco"
    $ws.Cells.Item($row, 3).Value  = "La Araucanía"
    $ws.Cells.Item($row, 4).Value  = $r.D
    $ws.Cells.Item($row, 5).Value  = 9
    $ws.Cells.Item($row, 6).Value  = 100112004
    $ws.Cells.Item($row, 7).Value  = "Cebolla"
    $ws.Cells.Item($row, 8).Value  = $r.H
    $ws.Cells.Item($row, 9).Value  = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}

Write-Output "Inserted 5 new rows (811-815); sheet now spans to row 827."
